$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 7.25
$ws.Range("J10").Value = "Easter"
$ws.Range("B11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = "Easter"
$ws.Range("B12").Value = 0
$ws.Range("J12").Value = "Easter"

$ws.Range("N10").Select()
